$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: the single run
#      "产品名称：神秘香料高级柴茶"
#    becomes two runs:
#      bold run   "产品名称"
#      normal run "：神秘香料臻品印度奶茶"
#    We rebuild the whole original run via InsertXML so every rPr flag -
#    including the <w:bCs/> on the new bold run - matches exactly, because
#    Font.Bold/Font.BoldBi manipulation on sub-ranges does not reliably
#    reproduce isolated <w:bCs/> state per run in this host.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1Full = $d.Range($p1.Start, $p1.Start + 13)

$splitXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>产品名称</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>：神秘香料臻品印度奶茶</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$p1Full.InsertXML($splitXml)

# ---------------------------------------------------------------------------
# 2) "主要特点：" heading becomes bold and its text changes to "主要功能："
# ---------------------------------------------------------------------------
$headingRange = $d.Content.Duplicate
$headingRange.Find.Execute("主要特点：", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingRange.Font.Bold = 1

$d.Content.Find.Execute("主要特点：", $true, $false, $false, $false, $false, $true, 1, $false, "主要功能：", 2)

# ---------------------------------------------------------------------------
# 3) Remaining plain text swaps (bullet titles + bodies)
# ---------------------------------------------------------------------------
$replacements = @(
    @("正宗混合", "正宗配方"),
    @("：我们的柴是优质黑茶叶的和谐混合，也是地香料的标志性选择，包括肉桂、豆瓜、丁香、姜和黑胡椒。", "：我们的奶茶选用优质黑茶，与肉桂、豆蔻、丁香、姜和黑胡椒等多种特色香料完美融合。"),
    @("健康增强成分", "成分更加健康"),
    @("：神秘香料柴茶中的每个成分都是出于自然健康益处而选择的。", "：神秘香料奶茶臻选自然原料，有利于健康。"),
    @("浓郁的香气和味道", "香气浓郁、口味醇厚"),
    @("：温暖，辣味和深，令人振奋的味道，我们的柴使它成为完美的饮料，开始你的一天或放松在晚上。", "：我们的奶茶气味温辛、口感醇厚，提神醒脑，是开启美好一天或晚上放松身心的完美饮品。"),
    @("多才多艺的酿造选项", "多元化的烹制选项"),
    @("：无论你喜欢你的柴热，作为一个令人耳目一新的冰茶，或作为奶油拿铁，我们的混合是多才多艺的，以满足任何偏好。", "：无论你是喜欢温热的奶茶，还是令人耳目一新的冰茶，或者是奶油拿铁，这款产品可以满足任何偏好。"),
    @("可持续来源", "原料可持续"),
    @("：致力于可持续性，我们从小型农场采购我们的成分，实践有机农业，不仅确保最好的品质，而且确保我们星球的福利。", "：我们注重可持续性，从小型农场采购原料，坚持有机农业，不仅能够确保极佳品质，而且可以确保对我们的星球有益。"),
    @("优雅的包装", "包装精致"),
    @("：神秘的香料柴茶是设计精美的生态友好包装，使其成为茶爱好者的理想礼物或豪华的礼物为自己。", "：神秘香料印度奶茶设计精美，采用生态友好的包装方式，因此是送给茶叶爱好者的理想礼物，也是送给自己的奢华之选。"),
    @("：我们站在产品后面，提供满意保证。", "：我们为产品背书，提供令人满意的保证。"),
    @("理想的选择", "适用人群"),
    @("：茶爱好者、有健康意识的个人、温暖、辛辣的饮料爱好者，以及任何希望探索传统印度柴的丰富口味的人。", "：茶叶爱好者、注重健康的个人、喜欢温辛饮料的群体，以及希望品尝传统印度奶茶丰富口感的人群。")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
